# Apply "hito 3" update to the HH / Presupuesto workbook.

$wb = $excel.ActiveWorkbook
$wsHH = $wb.Worksheets.Item("HH")
$wsPresupuesto = $wb.Worksheets.Item("Presupuesto")

# --- HH sheet ---------------------------------------------------------

# S3: update the hourly rate used in the formula (284.03 -> 314.25)
$wsHH.Range("S3").Formula = "=(314.25)*S2"

# Q6: extend the sum to include the new R12 entry
$wsHH.Range("Q6").Formula = "= Q8+R8+R9+R10+R11+R12"

# New entry R12 (hours logged), formatted like the other R-column entries ([h]:mm:ss)
$wsHH.Range("R12").Value = 1.2652777777777777
$wsHH.Range("R12").NumberFormat = "[h]:mm:ss"

# Q8 gets a new number format (date/time style), matching the new cellXf
$wsHH.Range("Q8").NumberFormat = "m/d/yy h:mm"

# Column Q needs to widen to fit the new date/time content (bestFit-style autosize)
$wsHH.Columns.Item(17).ColumnWidth = 13.45

# Update the selection on the HH sheet to S4
$wsHH.Range("S4").Select()

# --- Presupuesto sheet -------------------------------------------------

# Update the selection on the Presupuesto sheet to F10
$wsPresupuesto.Range("F10").Select()

$wb.Save()
